$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header "newest 2" in column G, row 1 (matches diff: new shared string
# "newest 2" and new cell G1 referencing it, with dimension expanding to A1:G7)
$ws.Range("G1").Value = "newest 2"
